# Updated simulated-game transition probability matrix on Sheet1
# (values recomputed after adding more simulated games; see commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state 0)
$ws.Range("B2").Value = 0.1413612565445026
$ws.Range("C2").Value = 0.6701570680628273
$ws.Range("J2").Value = 0.005235602094240838
$ws.Range("P2").Value = 0.1256544502617801
$ws.Range("S2").Value = 0.05759162303664921

# Row 3 (state 1)
$ws.Range("B3").Value = 0.02290076335877863
$ws.Range("C3").Value = 0.03053435114503817
$ws.Range("P3").Value = 0.8091603053435115
$ws.Range("S3").Value = 0.1374045801526718

# Row 4 (state 2)
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.7560975609756098
$ws.Range("S4").Value = 0.1707317073170732

# Row 5 (state 3)
$ws.Range("P5").Value = 1

# Row 6 (state 4)
$ws.Range("B6").Value = 0.04926108374384237
$ws.Range("D6").Value = 0.004926108374384237
$ws.Range("F6").Value = 0.05911330049261083
$ws.Range("J6").Value = 0.2906403940886699
$ws.Range("Q6").Value = 0.1921182266009852
$ws.Range("R6").Value = 0.07389162561576355
$ws.Range("S6").Value = 0.3300492610837438

# Row 7 (state 5)
$ws.Range("B7").Value = 0.06422018348623854
$ws.Range("D7").Value = 0.04587155963302753
$ws.Range("F7").Value = 0.1009174311926606
$ws.Range("J7").Value = 0.09174311926605505
$ws.Range("O7").Value = 0.03669724770642202
$ws.Range("Q7").Value = 0.1284403669724771
$ws.Range("R7").Value = 0.1926605504587156
$ws.Range("S7").Value = 0.3394495412844037

# Row 8 (state 6)
$ws.Range("B8").Value = 0.04775280898876404
$ws.Range("D8").Value = 0.03370786516853932
$ws.Range("F8").Value = 0.08707865168539326
$ws.Range("J8").Value = 0.09269662921348315
$ws.Range("O8").Value = 0.01123595505617977
$ws.Range("Q8").Value = 0.1769662921348314
$ws.Range("R8").Value = 0.1544943820224719
$ws.Range("S8").Value = 0.3960674157303371

# Row 9 (state 7)
$ws.Range("B9").Value = 0.09333333333333334
$ws.Range("D9").Value = 0.01333333333333333
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.04666666666666667
$ws.Range("O9").Value = 0.02666666666666667
$ws.Range("Q9").Value = 0.14
$ws.Range("R9").Value = 0.1466666666666667
$ws.Range("S9").Value = 0.4533333333333333

# Row 10 (state 8)
$ws.Range("B10").Value = 0.1073558648111332
$ws.Range("D10").Value = 0.02087475149105368
$ws.Range("E10").Value = 0.0009940357852882703
$ws.Range("F10").Value = 0.06858846918489066
$ws.Range("J10").Value = 0.1242544731610338
$ws.Range("O10").Value = 0.009940357852882704
$ws.Range("Q10").Value = 0.1848906560636183
$ws.Range("R10").Value = 0.1153081510934394
$ws.Range("S10").Value = 0.36779324055666

# Row 11 (state 9)
$ws.Range("F11").Value = 0.005952380952380952
$ws.Range("G11").Value = 0.1130952380952381
$ws.Range("J11").Value = 0.130952380952381
$ws.Range("K11").Value = 0.1785714285714286
$ws.Range("L11").Value = 0.5714285714285714

# Row 12 (state 10)
$ws.Range("G12").Value = 0.7669902912621359
$ws.Range("J12").Value = 0.145631067961165
$ws.Range("K12").Value = 0.009708737864077669
$ws.Range("L12").Value = 0.06796116504854369
$ws.Range("S12").Value = 0.009708737864077669

# Row 13 (state 11)
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333

# Row 15 (state 13)
$ws.Range("F15").Value = 0.03012048192771084
$ws.Range("H15").Value = 0.2349397590361446
$ws.Range("I15").Value = 0.0783132530120482
$ws.Range("J15").Value = 0.3493975903614458
$ws.Range("K15").Value = 0.04819277108433735
$ws.Range("M15").Value = 0.006024096385542169
$ws.Range("O15").Value = 0.06626506024096386
$ws.Range("S15").Value = 0.1867469879518072

# Row 16 (state 14)
$ws.Range("F16").Value = 0.05194805194805195
$ws.Range("H16").Value = 0.1558441558441558
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4545454545454545
$ws.Range("K16").Value = 0.1103896103896104
$ws.Range("M16").Value = 0.006493506493506494
$ws.Range("N16").Value = 0.006493506493506494
$ws.Range("O16").Value = 0.03246753246753246
$ws.Range("S16").Value = 0.1103896103896104

# Row 17 (state 15)
$ws.Range("F17").Value = 0.028125
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = 0.096875
$ws.Range("J17").Value = 0.434375
$ws.Range("K17").Value = 0.075
$ws.Range("M17").Value = 0.015625
$ws.Range("O17").Value = 0.05625
$ws.Range("S17").Value = 0.09375

# Row 18 (state 16)
$ws.Range("F18").Value = 0.03478260869565217
$ws.Range("H18").Value = 0.2043478260869565
$ws.Range("I18").Value = 0.1260869565217391
$ws.Range("J18").Value = 0.3565217391304348
$ws.Range("K18").Value = 0.07391304347826087
$ws.Range("M18").Value = 0.008695652173913044
$ws.Range("O18").Value = 0.08695652173913043
$ws.Range("S18").Value = 0.108695652173913

# Row 19 (state 17)
$ws.Range("F19").Value = 0.02011173184357542
$ws.Range("H19").Value = 0.2022346368715084
$ws.Range("I19").Value = 0.07150837988826815
$ws.Range("J19").Value = 0.4301675977653631
$ws.Range("K19").Value = 0.07709497206703911
$ws.Range("M19").Value = 0.01787709497206704
$ws.Range("N19").Value = 0.00111731843575419
$ws.Range("O19").Value = 0.07486033519553073
$ws.Range("S19").Value = 0.1050279329608939
